{"js": "// Applies the highlight-markup edits described in the commit diff:\n//   1) In the \"2. El (la) estudiante ...\" paragraph, highlight \"c\u00e9dula de\n//      ciudadan\u00eda #\" and \"Programa d\" in yellow.\n//   2) In the \"CODIGO Y PLAN ------\" table cell paragraph, split off the\n//      trailing \"------\" into its own (still bold) run and highlight it,\n//      along with the rest of that paragraph's text, in yellow.\n//   3) Replace \"Olga C.\" / \"Carolina Ch.\" with a yellow-highlighted \"\u2026\"\n//      after \"Proyect\u00f3: \" / \"Revis\u00f3: \" respectively.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Helper: find the first paragraph whose text satisfies predicate.\nfunction findParagraph(pred) {\n  for (const p of paragraphs.items) {\n    if (pred(p.text)) {\n      return p;\n    }\n  }\n  return null;\n}\n\n// ---------------------------------------------------------------------\n// 1) Paragraph \"2. El (la) estudiante [NOMBRE_ESTUDIANTE], ...\"\n// ---------------------------------------------------------------------\nconst estudianteParagraph = findParagraph(\n  (t) => t.indexOf(\"2. El (la) estudiante\") === 0\n);\n\nif (estudianteParagraph) {\n  const cedulaResults = estudianteParagraph.search(\"c\u00e9dula de ciudadan\u00eda #\", {\n    matchCase: false,\n  });\n  cedulaResults.load(\"items\");\n  const programaResults = estudianteParagraph.search(\"Programa d\", {\n    matchCase: true,\n  });\n  programaResults.load(\"items\");\n  await context.sync();\n\n  if (cedulaResults.items.length > 0) {\n    cedulaResults.items[0].font.highlightColor = \"yellow\";\n  }\n  if (programaResults.items.length > 0) {\n    programaResults.items[0].font.highlightColor = \"yellow\";\n  }\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// 2) Table cell paragraph \"CODIGO Y PLAN ------ Asignatura a homologar...\"\n// ---------------------------------------------------------------------\nconst codigoPlanParagraph = findParagraph(\n  (t) => t.indexOf(\"CODIGO Y PLAN\") === 0\n);\n\nif (codigoPlanParagraph) {\n  const dashResults = codigoPlanParagraph.search(\"------\", {\n    matchCase: true,\n  });\n  dashResults.load(\"items\");\n  const asignaturaResults = codigoPlanParagraph.search(\n    \"Asignatura a homologar en el \",\n    { matchCase: true }\n  );\n  asignaturaResults.load(\"items\");\n  const periodoResults = codigoPlanParagraph.search(\n    \"----- Periodo de 2025\",\n    { matchCase: true }\n  );\n  periodoResults.load(\"items\");\n  await context.sync();\n\n  if (\n    dashResults.items.length > 0 &&\n    asignaturaResults.items.length > 0 &&\n    periodoResults.items.length > 0\n  ) {\n    const dashRun = dashResults.items[0];\n    const asignaturaRun = asignaturaResults.items[0];\n    const periodoRun = periodoResults.items[0];\n\n    // Highlight \"------\" (splits it off from \"CODIGO Y PLAN \" automatically).\n    dashRun.font.highlightColor = \"yellow\";\n\n    // Isolate and highlight the single space run that sits between\n    // \"------\" and \"Asignatura a homologar en el \".\n    const afterDash = dashRun.getRange(\"After\");\n    const beforeAsignatura = asignaturaRun.getRange(\"Start\");\n    const spaceBetween = afterDash.expandTo(beforeAsignatura);\n    spaceBetween.font.highlightColor = \"yellow\";\n\n    asignaturaRun.font.highlightColor = \"yellow\";\n    periodoRun.font.highlightColor = \"yellow\";\n\n    await context.sync();\n  }\n}\n\n// ---------------------------------------------------------------------\n// 3) \"Proyect\u00f3: Olga C.\" / \"Revis\u00f3: Carolina Ch.\" paragraphs\n// ---------------------------------------------------------------------\nconst proyectoParagraph = findParagraph(\n  (t) => t.indexOf(\"Proyect\u00f3:\") === 0\n);\nconst revisoParagraph = findParagraph((t) => t.indexOf(\"Revis\u00f3:\") === 0);\n\nasync function replaceSignerWithEllipsis(paragraph, signerText) {\n  if (!paragraph) return;\n  const nameResults = paragraph.search(signerText, { matchCase: true });\n  nameResults.load(\"items\");\n  await context.sync();\n  if (nameResults.items.length === 0) return;\n\n  nameResults.items[0].insertText(\"\u2026\", \"Replace\");\n  await context.sync();\n\n  const ellipsisResults = paragraph.search(\"\u2026\", { matchCase: true });\n  ellipsisResults.load(\"items\");\n  await context.sync();\n  if (ellipsisResults.items.length > 0) {\n    ellipsisResults.items[0].font.highlightColor = \"yellow\";\n    await context.sync();\n  }\n}\n\nawait replaceSignerWithEllipsis(proyectoParagraph, \"Olga C.\");\nawait replaceSignerWithEllipsis(revisoParagraph, \"Carolina Ch.\");\n", "ps1": "# Applies the highlight-markup edits described in the commit diff:\n#   1) In the \"2. El (la) estudiante ...\" paragraph, highlight \"c\u00e9dula de\n#      ciudadan\u00eda #\" and \"Programa d\" in yellow.\n#   2) In the \"CODIGO Y PLAN ------\" table cell paragraph, split off the\n#      trailing \"------\" into its own (still bold) run and highlight it,\n#      along with the rest of that paragraph's text, in yellow.\n#   3) Replace \"Olga C.\" / \"Carolina Ch.\" with a yellow-highlighted \"\u2026\"\n#      after \"Proyect\u00f3: \" / \"Revis\u00f3: \" respectively.\n\n$d = $word.ActiveDocument\n$wdYellowHighlight = 7\n\nfunction Find-ParagraphStartingWith($doc, [string]$prefix) {\n    foreach ($p in $doc.Paragraphs) {\n        if ($p.Range.Text.StartsWith($prefix)) {\n            return $p\n        }\n    }\n    return $null\n}\n\nfunction Highlight-TextInRange($range, [string]$text, [bool]$matchCase) {\n    # Runs Find scoped to $range and, if found, highlights the matched\n    # text in yellow via the Selection object (Range objects returned\n    # from Find inside nested scopes such as paragraphs/cells do not\n    # reliably persist direct property writes, but Selection does).\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $text\n    $find.MatchCase = $matchCase\n    $find.Forward = $true\n    $find.Wrap = 0\n    $found = $find.Execute()\n    if ($found) {\n        $range.Select()\n        $word.Selection.Font.HighlightColorIndex = $wdYellowHighlight\n    }\n    return $found\n}\n\n# ---------------------------------------------------------------------\n# 1) Paragraph \"2. El (la) estudiante [NOMBRE_ESTUDIANTE], ...\"\n# ---------------------------------------------------------------------\n$estudianteParagraph = Find-ParagraphStartingWith $d \"2. El (la) estudiante\"\n\nif ($estudianteParagraph -ne $null) {\n    $rngCedula = $estudianteParagraph.Range\n    Highlight-TextInRange $rngCedula \"c\u00e9dula de ciudadan\u00eda #\" $false | Out-Null\n\n    $rngPrograma = $estudianteParagraph.Range\n    Highlight-TextInRange $rngPrograma \"Programa d\" $true | Out-Null\n}\n\n# ---------------------------------------------------------------------\n# 2) Table cell paragraph \"CODIGO Y PLAN ------ Asignatura a homologar...\"\n# ---------------------------------------------------------------------\n$codigoPlanParagraph = Find-ParagraphStartingWith $d \"CODIGO Y PLAN\"\n\nif ($codigoPlanParagraph -ne $null) {\n    $cellRange = $codigoPlanParagraph.Range\n\n    # Highlight the trailing \"------\" (splits it off from \"CODIGO Y PLAN \").\n    $rngDash = $cellRange.Duplicate\n    $foundDash = Highlight-TextInRange $rngDash \"------\" $true\n    $dashEnd = $rngDash.End\n\n    # Highlight \"Asignatura a homologar en el \".\n    $rngAsig = $cellRange.Duplicate\n    $foundAsig = Highlight-TextInRange $rngAsig \"Asignatura a homologar en el \" $true\n    $asigStart = $rngAsig.Start\n\n    # Highlight the single space run between \"------\" and \"Asignatura...\".\n    if ($foundDash -and $foundAsig) {\n        $rngSpace = $d.Range($dashEnd, $asigStart)\n        $rngSpace.Select()\n        $word.Selection.Font.HighlightColorIndex = $wdYellowHighlight\n    }\n\n    # Highlight \"----- Periodo de 2025\".\n    $rngPeriodo = $cellRange.Duplicate\n    Highlight-TextInRange $rngPeriodo \"----- Periodo de 2025\" $true | Out-Null\n}\n\n# ---------------------------------------------------------------------\n# 3) \"Proyect\u00f3: Olga C.\" / \"Revis\u00f3: Carolina Ch.\" paragraphs\n# ---------------------------------------------------------------------\nfunction Replace-SignerWithEllipsis($doc, [string]$paragraphPrefix, [string]$signerText) {\n    $paragraph = Find-ParagraphStartingWith $doc $paragraphPrefix\n    if ($paragraph -eq $null) {\n        return\n    }\n\n    $rngName = $paragraph.Range\n    $findName = $rngName.Find\n    $findName.ClearFormatting()\n    $findName.Text = $signerText\n    $findName.MatchCase = $true\n    $findName.Forward = $true\n    $findName.Wrap = 0\n    if ($findName.Execute()) {\n        $rngName.Text = \"\u2026\"\n    }\n\n    $rngEllipsis = $paragraph.Range\n    Highlight-TextInRange $rngEllipsis \"\u2026\" $true | Out-Null\n}\n\nReplace-SignerWithEllipsis $d \"Proyect\u00f3:\" \"Olga C.\"\nReplace-SignerWithEllipsis $d \"Revis\u00f3:\" \"Carolina Ch.\"\n"}
